# Dkk2-Lrp5.xlsx was regenerated with an updated TPM matrix; the ligand/receptor
# detection-rate, expression, and edge-specificity statistics for every
# Sending/Target cluster pair need to be refreshed to the newly computed values.
# (The set of clusters and their cell text is unchanged - only the numeric
# NATMI output columns E..J, M..T are recalculated.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 2: ECs -> ECs -------------------------------------------------
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3447283333333333
$ws.Range("H2").Value = 1.034185
$ws.Range("I2").Value = 0.07368549602308437
$ws.Range("J2").Value = 0.07368549602308436
$ws.Range("M2").Value = 17.08155333333333
$ws.Range("N2").Value = 51.24466
$ws.Range("O2").Value = 0.3501540759902865
$ws.Range("P2").Value = 0.3501540759902865
$ws.Range("Q2").Value = 5.888495411344444
$ws.Range("R2").Value = 52.99645870209999
$ws.Range("S2").Value = 0.02580127677384903
$ws.Range("T2").Value = 0.02580127677384903

# --- Row 3: ECs -> FAPs -----------------------------
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3447283333333333
$ws.Range("H3").Value = 1.034185
$ws.Range("I3").Value = 0.07368549602308437
$ws.Range("J3").Value = 0.07368549602308436
$ws.Range("O3").Value = 0.2142771237573249
$ws.Range("P3").Value = 0.2142771237573249
$ws.Range("Q3").Value = 3.603470433501666
$ws.Range("R3").Value = 32.431233901515
$ws.Range("S3").Value = 0.01578911615045832
$ws.Range("T3").Value = 0.01578911615045832

# --- Row 4: ECs -> Inflammatory-Mac -----------------
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3447283333333333
$ws.Range("H4").Value = 1.034185
$ws.Range("I4").Value = 0.07368549602308437
$ws.Range("J4").Value = 0.07368549602308436
$ws.Range("M4").Value = 8.398122666666666
$ws.Range("N4").Value = 25.194368
$ws.Range("O4").Value = 0.1721527793764119
$ws.Range("P4").Value = 0.1721527793764119
$ws.Range("Q4").Value = 2.895070830008888
$ws.Range("R4").Value = 26.05563747007999
$ws.Range("S4").Value = 0.01268516294010352
$ws.Range("T4").Value = 0.01268516294010352

# --- Row 5: ECs -> MuSCs ----------------------------
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3447283333333333
$ws.Range("H5").Value = 1.034185
$ws.Range("I5").Value = 0.07368549602308437
$ws.Range("J5").Value = 0.07368549602308436
$ws.Range("M5").Value = 4.514486333333333
$ws.Range("N5").Value = 13.543459
$ws.Range("O5").Value = 0.09254227409953211
$ws.Range("P5").Value = 0.09254227409953213
$ws.Range("Q5").Value = 1.556271349546111
$ws.Range("R5").Value = 14.006442145915
$ws.Range("S5").Value = 0.006819023370128257
$ws.Range("T5").Value = 0.006819023370128257

# --- Row 6: ECs -> Resolving-Mac --------------------
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3447283333333333
$ws.Range("H6").Value = 1.034185
$ws.Range("I6").Value = 0.07368549602308437
$ws.Range("J6").Value = 0.07368549602308436
$ws.Range("M6").Value = 8.335727666666667
$ws.Range("N6").Value = 25.007183
$ws.Range("O6").Value = 0.1708737467764446
$ws.Range("P6").Value = 0.1708737467764446
$ws.Range("Q6").Value = 2.873561505650555
$ws.Range("R6").Value = 25.862053550855
$ws.Range("S6").Value = 0.01259091678854524
$ws.Range("T6").Value = 0.01259091678854523

# --- Row 7: FAPs -> ECs -----------------------------
$ws.Range("I7").Value = 0.9050707286284559
$ws.Range("J7").Value = 0.9050707286284558
$ws.Range("M7").Value = 17.08155333333333
$ws.Range("N7").Value = 51.24466
$ws.Range("O7").Value = 0.3501540759902865
$ws.Range("P7").Value = 0.3501540759902865
$ws.Range("Q7").Value = 72.32773232334888
$ws.Range("R7").Value = 650.9495909101399
$ws.Range("S7").Value = 0.3169142046887523
$ws.Range("T7").Value = 0.3169142046887523

# --- Row 8: FAPs -> FAPs ----------------------------
$ws.Range("I8").Value = 0.9050707286284559
$ws.Range("J8").Value = 0.9050707286284558
$ws.Range("O8").Value = 0.2142771237573249
$ws.Range("P8").Value = 0.2142771237573249
$ws.Range("S8").Value = 0.1939359525274519
$ws.Range("T8").Value = 0.1939359525274518

# --- Row 9: FAPs -> Inflammatory-Mac ----------------
$ws.Range("I9").Value = 0.9050707286284559
$ws.Range("J9").Value = 0.9050707286284558
$ws.Range("M9").Value = 8.398122666666666
$ws.Range("N9").Value = 25.194368
$ws.Range("O9").Value = 0.1721527793764119
$ws.Range("P9").Value = 0.1721527793764119
$ws.Range("Q9").Value = 35.55983208318577
$ws.Range("R9").Value = 320.038488748672
$ws.Range("S9").Value = 0.155810441465623
$ws.Range("T9").Value = 0.155810441465623

# --- Row 10: FAPs -> MuSCs ---------------------------
$ws.Range("I10").Value = 0.9050707286284559
$ws.Range("J10").Value = 0.9050707286284558
$ws.Range("M10").Value = 4.514486333333333
$ws.Range("N10").Value = 13.543459
$ws.Range("O10").Value = 0.09254227409953211
$ws.Range("P10").Value = 0.09254227409953213
$ws.Range("Q10").Value = 19.11550739695122
$ws.Range("R10").Value = 172.039566572561
$ws.Range("S10").Value = 0.08375730344819782
$ws.Range("T10").Value = 0.08375730344819782

# --- Row 11: FAPs -> Resolving-Mac -------------------
$ws.Range("I11").Value = 0.9050707286284559
$ws.Range("J11").Value = 0.9050707286284558
$ws.Range("M11").Value = 8.335727666666667
$ws.Range("N11").Value = 25.007183
$ws.Range("O11").Value = 0.1708737467764446
$ws.Range("P11").Value = 0.1708737467764446
$ws.Range("Q11").Value = 35.29563545128411
$ws.Range("R11").Value = 317.660719061557
$ws.Range("S11").Value = 0.154652826498431
$ws.Range("T11").Value = 0.154652826498431

# --- Row 12: MuSCs -> ECs ----------------------------
$ws.Range("G12").Value = 0.027522
$ws.Range("H12").Value = 0.082566
$ws.Range("I12").Value = 0.005882812712079546
$ws.Range("J12").Value = 0.005882812712079545
$ws.Range("M12").Value = 17.08155333333333
$ws.Range("N12").Value = 51.24466
$ws.Range("O12").Value = 0.3501540759902865
$ws.Range("P12").Value = 0.3501540759902865
$ws.Range("Q12").Value = 0.47011851084
$ws.Range("R12").Value = 4.23106659756
$ws.Range("S12").Value = 0.002059890849422125
$ws.Range("T12").Value = 0.002059890849422124

# --- Row 13: MuSCs -> FAPs ---------------------------
$ws.Range("G13").Value = 0.027522
$ws.Range("H13").Value = 0.082566
$ws.Range("I13").Value = 0.005882812712079546
$ws.Range("J13").Value = 0.005882812712079545
$ws.Range("O13").Value = 0.2142771237573249
$ws.Range("P13").Value = 0.2142771237573249
$ws.Range("Q13").Value = 0.287689475106
$ws.Range("R13").Value = 2.589205275954
$ws.Range("S13").Value = 0.001260552187547433
$ws.Range("T13").Value = 0.001260552187547433

# --- Row 14: MuSCs -> Inflammatory-Mac ---------------
$ws.Range("G14").Value = 0.027522
$ws.Range("H14").Value = 0.082566
$ws.Range("I14").Value = 0.005882812712079546
$ws.Range("J14").Value = 0.005882812712079545
$ws.Range("M14").Value = 8.398122666666666
$ws.Range("N14").Value = 25.194368
$ws.Range("O14").Value = 0.1721527793764119
$ws.Range("P14").Value = 0.1721527793764119
$ws.Range("Q14").Value = 0.231133132032
$ws.Range("R14").Value = 2.080198188288
$ws.Range("S14").Value = 0.001012742558935382
$ws.Range("T14").Value = 0.001012742558935381

# --- Row 15: MuSCs -> MuSCs --------------------------
$ws.Range("G15").Value = 0.027522
$ws.Range("H15").Value = 0.082566
$ws.Range("I15").Value = 0.005882812712079546
$ws.Range("J15").Value = 0.005882812712079545
$ws.Range("M15").Value = 4.514486333333333
$ws.Range("N15").Value = 13.543459
$ws.Range("O15").Value = 0.09254227409953211
$ws.Range("P15").Value = 0.09254227409953213
$ws.Range("Q15").Value = 0.124247692866
$ws.Range("R15").Value = 1.118229235794
$ws.Range("S15").Value = 0.0005444088664774772
$ws.Range("T15").Value = 0.0005444088664774772

# --- Row 16: MuSCs -> Resolving-Mac ------------------
$ws.Range("G16").Value = 0.027522
$ws.Range("H16").Value = 0.082566
$ws.Range("I16").Value = 0.005882812712079546
$ws.Range("J16").Value = 0.005882812712079545
$ws.Range("M16").Value = 8.335727666666667
$ws.Range("N16").Value = 25.007183
$ws.Range("O16").Value = 0.1708737467764446
$ws.Range("P16").Value = 0.1708737467764446
$ws.Range("Q16").Value = 0.229415896842
$ws.Range("R16").Value = 2.064743071578
$ws.Range("S16").Value = 0.00100521824969713
$ws.Range("T16").Value = 0.00100521824969713

# --- Row 17: Resolving-Mac -> ECs --------------------
$ws.Range("G17").Value = 0.07186433333333334
$ws.Range("H17").Value = 0.215593
$ws.Range("I17").Value = 0.01536096263638017
$ws.Range("J17").Value = 0.01536096263638017
$ws.Range("M17").Value = 17.08155333333333
$ws.Range("N17").Value = 51.24466
$ws.Range("O17").Value = 0.3501540759902865
$ws.Range("P17").Value = 0.3501540759902865
$ws.Range("Q17").Value = 1.227554442597778
$ws.Range("R17").Value = 11.04798998338
$ws.Range("S17").Value = 0.005378703678263015
$ws.Range("T17").Value = 0.005378703678263015

# --- Row 18: Resolving-Mac -> FAPs -------------------
$ws.Range("G18").Value = 0.07186433333333334
$ws.Range("H18").Value = 0.215593
$ws.Range("I18").Value = 0.01536096263638017
$ws.Range("J18").Value = 0.01536096263638017
$ws.Range("O18").Value = 0.2142771237573249
$ws.Range("P18").Value = 0.2142771237573249
$ws.Range("Q18").Value = 0.7512031224296667
$ws.Range("R18").Value = 6.760828101867
$ws.Range("S18").Value = 0.003291502891867278
$ws.Range("T18").Value = 0.003291502891867278

# --- Row 19: Resolving-Mac -> Inflammatory-Mac -------
$ws.Range("G19").Value = 0.07186433333333334
$ws.Range("H19").Value = 0.215593
$ws.Range("I19").Value = 0.01536096263638017
$ws.Range("J19").Value = 0.01536096263638017
$ws.Range("M19").Value = 8.398122666666666
$ws.Range("N19").Value = 25.194368
$ws.Range("O19").Value = 0.1721527793764119
$ws.Range("P19").Value = 0.1721527793764119
$ws.Range("Q19").Value = 0.6035254866915555
$ws.Range("R19").Value = 5.431729380224
$ws.Range("S19").Value = 0.002644432411750063
$ws.Range("T19").Value = 0.002644432411750063

# --- Row 20: Resolving-Mac -> MuSCs ------------------
$ws.Range("G20").Value = 0.07186433333333334
$ws.Range("H20").Value = 0.215593
$ws.Range("I20").Value = 0.01536096263638017
$ws.Range("J20").Value = 0.01536096263638017
$ws.Range("M20").Value = 4.514486333333333
$ws.Range("N20").Value = 13.543459
$ws.Range("O20").Value = 0.09254227409953211
$ws.Range("P20").Value = 0.09254227409953213
$ws.Range("Q20").Value = 0.3244305506874444
$ws.Range("R20").Value = 2.919874956187
$ws.Range("S20").Value = 0.001421538414728565
$ws.Range("T20").Value = 0.001421538414728566

# --- Row 21: Resolving-Mac -> Resolving-Mac ----------
$ws.Range("G21").Value = 0.07186433333333334
$ws.Range("H21").Value = 0.215593
$ws.Range("I21").Value = 0.01536096263638017
$ws.Range("J21").Value = 0.01536096263638017
$ws.Range("M21").Value = 8.335727666666667
$ws.Range("N21").Value = 25.007183
$ws.Range("O21").Value = 0.1708737467764446
$ws.Range("P21").Value = 0.1708737467764446
$ws.Range("Q21").Value = 0.5990415116132223
$ws.Range("R21").Value = 5.391373604519001
$ws.Range("S21").Value = 0.002624785239771253
$ws.Range("T21").Value = 0.002624785239771253
